$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 49, shifting rows 49:87 down to 50:88
$ws.Rows.Item(49).Insert()

# Fill in the new row 49 with the new record
$ws.Cells.Item(49, 1).Value = 11
$ws.Cells.Item(49, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(49, 3).Value = "Bíobío"
$ws.Cells.Item(49, 4).NumberFormat = $ws.Cells.Item(50, 4).NumberFormat
$ws.Cells.Item(49, 4).Value = 44944
$ws.Cells.Item(49, 5).Value = 8
$ws.Cells.Item(49, 6).Value = 100112037
$ws.Cells.Item(49, 7).Value = "Cebollín"
$ws.Cells.Item(49, 8).Value = "Sin especificar"
$ws.Cells.Item(49, 9).Value = "Primera"
$ws.Cells.Item(49, 10).Value = 220
$ws.Cells.Item(49, 11).Value = 3000
$ws.Cells.Item(49, 12).Value = 3200
$ws.Cells.Item(49, 13).Value = 3091
$ws.Cells.Item(49, 14).Value = "`$/paquete 36 unidades"
$ws.Cells.Item(49, 15).Value = "Región Metropolitana"
$ws.Cells.Item(49, 16).Value = 86
$ws.Cells.Item(49, 17).Value = 36
$ws.Cells.Item(49, 18).Value = "Hortaliza"
